$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 326, shifting existing rows 326-406 down to 327-407
$ws.Rows.Item(326).Insert()

# Populate the new row 326 with the new data point
$ws.Cells.Item(326, 1).Value = 6
$ws.Cells.Item(326, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(326, 3).Value = "Metropolitana"
$ws.Cells.Item(326, 4).Value = 44551
$ws.Cells.Item(326, 5).Value = 13
$ws.Cells.Item(326, 6).Value = 100112044
$ws.Cells.Item(326, 7).Value = "Perejil"
$ws.Cells.Item(326, 8).Value = "Sin especificar"
$ws.Cells.Item(326, 9).Value = "Primera"
$ws.Cells.Item(326, 10).Value = 230
$ws.Cells.Item(326, 11).Value = 9000
$ws.Cells.Item(326, 12).Value = 10000
$ws.Cells.Item(326, 13).Value = 9348
$ws.Cells.Item(326, 14).Value = "$/docena de atados"
$ws.Cells.Item(326, 15).Value = "Región Metropolitana"
$ws.Cells.Item(326, 16).Value = 3116
$ws.Cells.Item(326, 17).Value = 3
$ws.Cells.Item(326, 18).Value = "Hortaliza"
